$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 179, shifting existing rows 179-182 down to 180-183
$ws.Rows.Item(179).Insert()

# Populate new row 179 with data
$ws.Range("A179").Value = 5
$ws.Range("B179").Value = "Macroferia Regional de Talca"
$ws.Range("C179").Value = "Maule"
$ws.Range("D179").Value = 44448
$ws.Range("D179").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E179").Value = 7
$ws.Range("F179").Value = 100114013
$ws.Range("G179").Value = "Zanahoria"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 400
$ws.Range("K179").Value = 5500
$ws.Range("L179").Value = 5500
$ws.Range("M179").Value = 5500
$ws.Range("N179").Value = "$/saco 20 kilos"
$ws.Range("O179").Value = "Región de Ñuble"
$ws.Range("P179").Value = 275
$ws.Range("Q179").Value = 20
$ws.Range("R179").Value = "Hortaliza"
